# Natmi following Dr Hou advice
# Update the LR-pair sheet so the FAPs/ECs cluster combinations are all
# represented (FAPs and ECs as both sender and target), replacing the old
# two-row summary with the full 6-row breakdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Wnt1","Fzd8","FAPs",1,0.3333333333333333,0.027123,0.081369,0.07131444737854614,0.07131444737854616,2,0.6666666666666666,1.027778333333333,3.083335,0.08020467841353289,0.08020467841353289,0.027876431735,0.250887885615,0.005719752318235107,0.005719752318235108),
    @("FAPs","Wnt1","Fzd8","ECs",1,0.3333333333333333,0.027123,0.081369,0.07131444737854614,0.07131444737854616,3,1,7.273511666666667,21.820535,0.567602609669802,0.567602609669802,0.197279456935,1.775515112415,0.04047826643922256,0.04047826643922257),
    @("FAPs","Wnt1","Fzd8","sCs",1,0.3333333333333333,0.027123,0.081369,0.07131444737854614,0.07131444737854616,3,1,4.513153666666667,13.539461,0.3521927119166651,0.3521927119166651,0.122410266901,1.101692402109,0.02511642862108848,0.02511642862108848),
    @("ECs","Wnt1","Fzd8","FAPs",3,1,0.3532066666666667,1.05962,0.9286855526214538,0.9286855526214538,2,0.6666666666666666,1.027778333333333,3.083335,0.08020467841353289,0.08020467841353289,0.3630181591888889,3.2671634327,0.07448492609529778,0.07448492609529778),
    @("ECs","Wnt1","Fzd8","ECs",3,1,0.3532066666666667,1.05962,0.9286855526214538,0.9286855526214538,3,1,7.273511666666667,21.820535,0.567602609669802,0.567602609669802,2.569052810744445,23.1214752967,0.5271243432305794,0.5271243432305794),
    @("ECs","Wnt1","Fzd8","sCs",3,1,0.3532066666666667,1.05962,0.9286855526214538,0.9286855526214538,3,1,4.513153666666667,13.539461,0.3521927119166651,0.3521927119166651,1.594075962757778,14.34668366482,0.3270762832955766,0.3270762832955766)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = 1 + $j
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
